$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 (channel 7023, delay 0) already holds the values we want in row 2 -
# copy it over (values only, so we don't drag any border/date formatting
# along with it) to replace the old 7002 / 2 entry in row 2.
$ws.Range("A23:B23").Copy()
$ws.Range("A2").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = $false
